{"js": "// Remove the \"Ver no Jupiter...\" / copyright footer block that follows the\n// \"LOQ4073: Qu\u00edmica Geral II (Requisito fraco)\" requirements line, along\n// with the blank paragraph that separates them.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the \"Requisitos\" line anchor paragraph by its known text.\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"LOQ4073\") !== -1) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex === -1) {\n  throw new Error(\"Could not find the 'LOQ4073' requirements paragraph.\");\n}\n\n// The three paragraphs immediately after the anchor are:\n//   1) a blank spacer paragraph\n//   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n//   3) the \"\u00a9 2020 ...\" footer line\n// All three are removed, leaving the anchor paragraph directly followed by\n// the (unrelated) trailing blank / page-break paragraphs.\nconst toDelete = [];\nfor (let i = anchorIndex + 1; i <= anchorIndex + 3; i++) {\n  if (i < items.length) {\n    toDelete.push(items[i]);\n  }\n}\n\ntoDelete.forEach((p) => p.delete());\nawait context.sync();\n", "ps1": "# Remove the \"Ver no Jupiter...\" / copyright footer block that follows the\n# \"LOQ4073: Qu\u00edmica Geral II (Requisito fraco)\" requirements line, along\n# with the blank paragraph that separates them.\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph (\"LOQ4073: ...\") via Find so the script does\n# not depend on a hard-coded paragraph number.\n$finder = $d.Content\n$finder.Find.Execute(\"LOQ4073\") | Out-Null\n$anchorIndex = $finder.Paragraphs.Item(1).Index\n\n# The three paragraphs immediately following the anchor are:\n#   1) a blank spacer paragraph\n#   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n#   3) the \"\u00a9 2020 ...\" footer line\n# Delete them (each deletion shifts later paragraphs up, so we repeatedly\n# remove the paragraph right after the anchor).\nfor ($n = 1; $n -le 3; $n++) {\n    $d.Paragraphs.Item($anchorIndex + 1).Range.Delete()\n}\n"}
